# Apply the crawl refresh update for bread_coop_2023-01-30.xlsx
# - Refresh the timestamp column (O) for every data row to the new crawl time
# - Update a couple of product rows whose scraped attributes changed between crawls

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-01-30 12:58:29"

$lastRow = $ws.UsedRange.Rows.Count

# Refresh timestamp (column O) on every data row (rows 2..lastRow)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Row 131: Leisi Kuchenteig rund ausgewallt - declaration icons gained vegan/vegetarian tags
$ws.Range("N131").Value = "['chilled', 'vegan', 'vegetarian']"

# Row 251: Pasquier Schokobroetchen 16 Stueck - now flagged as out of stock online
$ws.Range("M251").Value = "Pasquier Schokobr" + [char]0xF6 + "tchen 16 St" + [char]0xFC + "ck - Online kein Bestand 8.50 Schweizer Franken"

# Row 374: Naturaplan Bio Appenzeller Baerli Biber 3x62g - now on a 20% promo
$ws.Range("G374").NumberFormat = "@"
$ws.Range("G374").Value = "3.80"
$ws.Range("H374").Value = "2.04/100g"
$ws.Range("J374").NumberFormat = "@"
$ws.Range("J374").Value = "2.04"
$ws.Range("M374").Value = "Naturaplan Bio Appenzeller B" + [char]0xE4 + "rli Biber 3x62g 20% Aktion 3.80 Schweizer Franken statt 4.80 Schweizer Franken"

# Row 375: Naturaplan Bio Bischofberger Biber Herzli 6x32g - now on a 20% promo
$ws.Range("G375").NumberFormat = "@"
$ws.Range("G375").Value = "4.40"
$ws.Range("H375").Value = "2.29/100g"
$ws.Range("J375").NumberFormat = "@"
$ws.Range("J375").Value = "2.29"
$ws.Range("M375").Value = "Naturaplan Bio Bischofberger Biber Herzli 6x32g 20% Aktion 4.40 Schweizer Franken statt 5.50 Schweizer Franken"
